$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextCell "D2" "68.148.60"

# Row 3 - Ethereum
Set-TextCell "D3" "2.637.43"
$ws.Range("E3").Value = "  +0.83%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.02%  "

# Row 5 - BNB
Set-TextCell "D5" "596.78"
$ws.Range("E5").Value = "  +0.15%  "

# Row 6 - Solana
Set-TextCell "D6" "154.76"
$ws.Range("E6").Value = "  +1.11%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.03%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +0.32%  "

# Row 9 - LidoStakedEther
Set-TextCell "D9" "2.636.54"
$ws.Range("E9").Value = "  +0.79%  "

# Row 10 - Dogecoin
Set-TextCell "D10" "0.144"
$ws.Range("E10").Value = "  +8.28%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -0.83%  "

# Row 12 - Toncoin
$ws.Range("E12").Value = "  +0.22%  "

# Row 13 - Cardano
Set-TextCell "D13" "0.350"
$ws.Range("E13").Value = "  +1.05%  "

# Row 14 - was ShibaInu, now Avalanche
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextCell "D14" "27.86"
$ws.Range("E14").Value = "  +1.15%  "

# Row 15 - was Avalanche, now ShibaInu
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextCell "D15" "0.0000192"
$ws.Range("E15").Value = "  +2.74%  "

# Row 16 - WrappedliquidstakedEther2.0
Set-TextCell "D16" "3.119.00"

# Row 17 - WrappedBTC
Set-TextCell "D17" "68.043.83"
$ws.Range("E17").Value = "  +0.44%  "

# Row 18 - WrappedEther
Set-TextCell "D18" "2.631.68"
$ws.Range("E18").Value = "  +0.50%  "

# Row 19 - Chainlink
Set-TextCell "D19" "11.36"
$ws.Range("E19").Value = "  +0.28%  "

# Row 20 - BitcoinCash
Set-TextCell "D20" "362.73"
$ws.Range("E20").Value = "  -1.10%  "

# Row 21 - Uniswap
Set-TextCell "D21" "7.42"
$ws.Range("E21").Value = "  +0.55%  "

# Row 22 - Polkadot
Set-TextCell "D22" "4.36"
$ws.Range("E22").Value = "  +3.55%  "

# Row 23 - NEARProtocol
$ws.Range("E23").Value = "  +0.64%  "

# Row 24 - SuiNetwork
$ws.Range("E24").Value = "  -0.72%  "

# Row 25 - Litecoin
Set-TextCell "D25" "74.91"
$ws.Range("E25").Value = "  +3.01%  "

# Row 26 - Dai
Set-TextCell "D26" "0.999"
$ws.Range("E26").Value = "  -0.05%  "

# Row 27 - Aptos
Set-TextCell "D27" "9.70"
$ws.Range("E27").Value = "  -2.16%  "

# Row 28 - PEPE
$ws.Range("E28").Value = "  +1.68%  "

# Row 29 - WrappedeETH
Set-TextCell "D29" "2.769.50"
$ws.Range("E29").Value = "  +0.69%  "

# Row 30 - Binance-PegBSC-USD
$ws.Range("E30").Value = "  +0.24%  "

# Row 31 - Bittensor
Set-TextCell "D31" "560.98"
$ws.Range("E31").Value = "  -1.12%  "

# Row 32 - InternetComputer(DFINITY)
Set-TextCell "D32" "7.99"
$ws.Range("E32").Value = "  +1.73%  "

# Row 33 - Fetch.AI
$ws.Range("E33").Value = "  +0.41%  "

# Row 34 - PancakeSwap
$ws.Range("E34").Value = "  +1.28%  "

# Row 35 - Kaspa
$ws.Range("E35").Value = "  +2.53%  "

# Row 36 - FirstDigitalUSD
Set-TextCell "D36" "1.00"
$ws.Range("E36").Value = "  +0.11%  "

# Row 37 - ImmutableX
$ws.Range("E37").Value = "  +4.03%  "

# Row 38 - Monero
Set-TextCell "D38" "160.45"
$ws.Range("E38").Value = "  -0.76%  "

# Row 39 - EthereumClassic
Set-TextCell "D39" "19.29"
$ws.Range("E39").Value = "  +1.25%  "

# Row 40 - PolygonEcosystemToken
$ws.Range("E40").Value = "  +1.60%  "

# Row 41 - Stacks
$ws.Range("E41").Value = "  -0.23%  "

# Row 42 - RenderToken
$ws.Range("E42").Value = "  +0.53%  "

# Row 43 - BabyDogeCoin
$ws.Range("E43").Value = "  +3.73%  "

# Row 44 - WhiteBITCoin
$ws.Range("E44").Value = "  +2.54%  "

# Row 45 - dogwifhat
$ws.Range("E45").Value = "  -0.59%  "

# Row 46 - USDe
$ws.Range("E46").Value = "  +0.01%  "

# Row 47 - OKB
Set-TextCell "D47" "40.43"
$ws.Range("E47").Value = "  +0.78%  "

# Row 48 - Aave
Set-TextCell "D48" "158.16"
$ws.Range("E48").Value = "  +2.19%  "

# Row 49 - Filecoin
$ws.Range("E49").Value = "  +1.98%  "

# Row 50 - InjectiveProtocol
Set-TextCell "D50" "21.92"

# Row 51 - Cronos
$ws.Range("E51").Value = "  +1.48%  "
